$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row onto the two new rows
$ws.Range("A5:B5").Copy()
$ws.Range("A6:B7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add the two new rows of data
$ws.Range("A6").Value = "Companhia Siderúrgica Nacional"
$ws.Range("B6").Value = "CSNA3:BVMF"
$ws.Range("A7").Value = "Eletrobras"
$ws.Range("B7").Value = "ELET6:BVMF"

# Widen column A to fit the new, longer company name
$ws.Columns.Item(1).ColumnWidth = 37.29

# Update the active cell selection
$ws.Range("E9").Select()
